$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1710.9844
$ws.Range("I15").Value = 1710.9844
$ws.Range("K15").Value = 5132.9532
$ws.Range("M15").Value = -4963.9532
$ws.Range("H40").Value = 5000
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("I55").Value = 223.25
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 223.25
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -9.25
$ws.Range("N55").ClearContents()
$ws.Range("H69").Value = 4207.5
$ws.Range("I69").Value = 3500
$ws.Range("K69").Value = 10500
$ws.Range("M69").Value = -9626
$ws.Range("H70").Value = 1616.5
$ws.Range("I70").Value = 1049.75
$ws.Range("J70").Value = 2750
$ws.Range("K70").Value = 3149.25
$ws.Range("L70").Value = 8250
$ws.Range("M70").Value = -2879.25
$ws.Range("N70").Value = -8790
$ws.Range("H72").Value = 4207.5
$ws.Range("I72").Value = 3500
$ws.Range("K72").Value = 31500
$ws.Range("M72").Value = -27132
$ws.Range("H73").Value = 1616.5
$ws.Range("I73").Value = 1049.75
$ws.Range("J73").Value = 2750
$ws.Range("K73").Value = 3149.25
$ws.Range("L73").Value = 8250
$ws.Range("M73").Value = -2213.25
$ws.Range("N73").Value = -10122
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H112").Value = 2500
$ws.Range("J112").Value = 2500
$ws.Range("L112").Value = 7500
$ws.Range("N112").Value = -9716
$ws.Range("H115").Value = 10200
$ws.Range("J115").Value = 12000
$ws.Range("L115").Value = 36000
$ws.Range("N115").Value = -39134
$ws.Range("H116").Value = 4058.5293
$ws.Range("I116").Value = 3922.923
$ws.Range("K116").Value = 3922.923
$ws.Range("M116").Value = -480.9229999999998
$ws.Range("H137").Value = 2388.3489
$ws.Range("I137").Value = 1611.1786
$ws.Range("K137").Value = 4833.5358
$ws.Range("M137").Value = -2283.5358
$ws.Range("H138").Value = 5283.727
$ws.Range("J138").Value = 5283.727
$ws.Range("L138").Value = 15851.181
$ws.Range("N138").Value = -26131.181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16168.792
$ws.Range("I32").Value = 16168.792
$ws.Range("K32").Value = 16168.792
$ws.Range("M32").Value = -15881.792
$ws.Range("H74").Value = 2212.32
$ws.Range("I74").Value = 1383.8667
$ws.Range("K74").Value = 1383.8667
$ws.Range("M74").Value = -509.8667
$ws.Range("H77").Value = 2212.32
$ws.Range("I77").Value = 1383.8667
$ws.Range("K77").Value = 6919.333500000001
$ws.Range("M77").Value = -2551.333500000001
$ws.Range("H97").Value = 279.25
$ws.Range("I97").Value = 279.25
$ws.Range("K97").Value = 279.25
$ws.Range("M97").Value = 216.75
$ws.Range("H110").Value = 1039.4
$ws.Range("I110").Value = 1066.2222
$ws.Range("J110").Value = 798
$ws.Range("K110").Value = 1066.2222
$ws.Range("L110").Value = 798
$ws.Range("M110").Value = 978.7778000000001
$ws.Range("N110").Value = -4888
$ws.Range("H114").Value = 70000
$ws.Range("J114").Value = 70000
$ws.Range("L114").Value = 70000
$ws.Range("N114").Value = -78678

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1331.8334
$ws.Range("J64").Value = 1663.6666
$ws.Range("L64").Value = 1663.6666
$ws.Range("N64").Value = -2113.6666
$ws.Range("H67").Value = 1331.8334
$ws.Range("J67").Value = 1663.6666
$ws.Range("L67").Value = 1663.6666
$ws.Range("N67").Value = -3223.6666
$ws.Range("H86").Value = 4223.2144
$ws.Range("I86").Value = 4075.182
$ws.Range("J86").Value = 4766
$ws.Range("K86").Value = 4075.182
$ws.Range("L86").Value = 4766
$ws.Range("M86").Value = -2952.182
$ws.Range("N86").Value = -7012
$ws.Range("H89").Value = 4223.2144
$ws.Range("I89").Value = 4075.182
$ws.Range("J89").Value = 4766
$ws.Range("K89").Value = 20375.91
$ws.Range("L89").Value = 23830
$ws.Range("M89").Value = -14759.91
$ws.Range("N89").Value = -35062
$ws.Range("H134").Value = 4240.16
$ws.Range("I134").Value = 4087.1738
$ws.Range("K134").Value = 12261.5214
$ws.Range("M134").Value = -9726.5214

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H53").Value = 62342
$ws.Range("J53").Value = 62342
$ws.Range("L53").Value = 62342
$ws.Range("N53").Value = -63556
$ws.Range("H141").Value = 90399.28999999999
$ws.Range("J141").Value = 88799.164
$ws.Range("L141").Value = 88799.164
$ws.Range("N141").Value = -99159.164

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 2996.3333
$ws.Range("J92").Value = 2996.3333
$ws.Range("L92").Value = 8988.999899999999
$ws.Range("N92").Value = -11484.9999
$ws.Range("H121").Value = 1444
$ws.Range("I121").Value = 501.5
$ws.Range("J121").Value = 2700.6667
$ws.Range("K121").Value = 1504.5
$ws.Range("L121").Value = 8102.000100000001
$ws.Range("M121").Value = -194.5
$ws.Range("N121").Value = -10722.0001
$ws.Range("H122").Value = 102148
$ws.Range("J122").Value = 127310
$ws.Range("L122").Value = 1145790
$ws.Range("N122").Value = -1150690
$ws.Range("H131").Value = 2076
$ws.Range("I131").Value = 1400
$ws.Range("K131").Value = 4200
$ws.Range("M131").Value = 840

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2338
$ws.Range("I102").Value = 2000
$ws.Range("J102").Value = 3014
$ws.Range("K102").Value = 2000
$ws.Range("L102").Value = 3014
$ws.Range("M102").Value = -378
$ws.Range("N102").Value = -6258
$ws.Range("H122").Value = 1511.375
$ws.Range("I122").Value = 1511.375
$ws.Range("K122").Value = 4534.125
$ws.Range("M122").Value = -2084.125
$ws.Range("H132").Value = 2919.5
$ws.Range("I132").Value = 2239.6
$ws.Range("K132").Value = 6718.799999999999
$ws.Range("M132").Value = -4188.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H46").Value = 445
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H68").Value = 2586.875
$ws.Range("J68").Value = 2932.3333
$ws.Range("L68").Value = 2932.3333
$ws.Range("N68").Value = -4430.3333
$ws.Range("H71").Value = 2586.875
$ws.Range("J71").Value = 2932.3333
$ws.Range("L71").Value = 14661.6665
$ws.Range("N71").Value = -22149.6665
$ws.Range("H110").Value = 150000
$ws.Range("J110").Value = 150000
$ws.Range("L110").Value = 150000
$ws.Range("N110").Value = -158180

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8865.333000000001
$ws.Range("I81").Value = 5880.6665
$ws.Range("J81").Value = 14834.667
$ws.Range("K81").Value = 11761.333
$ws.Range("L81").Value = 29669.334
$ws.Range("M81").Value = -10700.333
$ws.Range("N81").Value = -31791.334
$ws.Range("H84").Value = 8865.333000000001
$ws.Range("I84").Value = 5880.6665
$ws.Range("J84").Value = 14834.667
$ws.Range("K84").Value = 58806.665
$ws.Range("L84").Value = 148346.67
$ws.Range("M84").Value = -53502.665
$ws.Range("N84").Value = -158954.67
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
